$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "mortgage"
$ws.Range("D2").Value = 0.476609855890274
$ws.Range("C3").Value = "lending"
$ws.Range("D3").Value = 0.4755072295665741
$ws.Range("C4").Value = "refinance"
$ws.Range("D4").Value = 0.4599137008190155
$ws.Range("C5").Value = "buyer"
$ws.Range("D5").Value = 0.4443101286888122
$ws.Range("C6").Value = "floor"
$ws.Range("D6").Value = 0.425670862197876
$ws.Range("C7").Value = "investor"
$ws.Range("D7").Value = 0.4170790910720825
$ws.Range("C8").Value = "tap"
$ws.Range("D8").Value = 0.3968437910079956
$ws.Range("C9").Value = "method"
$ws.Range("D9").Value = 0.3918479084968567
$ws.Range("C10").Value = "deposit"
$ws.Range("D10").Value = 0.3908049464225769
$ws.Range("D11").Value = 0.3821022808551788
$ws.Range("D12").Value = 0.5050563216209412
$ws.Range("D13").Value = 0.4787676632404327
$ws.Range("D14").Value = 0.4347586631774902
$ws.Range("C15").Value = "scattered"
$ws.Range("D15").Value = 0.4143691658973694
$ws.Range("C16").Value = "commodity"
$ws.Range("D16").Value = 0.3973327577114105
$ws.Range("C17").Value = "intensive"
$ws.Range("D17").Value = 0.3933558762073517
$ws.Range("C18").Value = "restaurateur"
$ws.Range("D18").Value = 0.3930320143699646
$ws.Range("C19").Value = "pressure"
$ws.Range("D19").Value = 0.3810980916023254
$ws.Range("C20").Value = "raise"
$ws.Range("D20").Value = 0.3708976805210113
$ws.Range("C21").Value = "respective"
$ws.Range("D21").Value = 0.3707828223705292
$ws.Range("D22").Value = 0.473964661359787
$ws.Range("C23").Value = "appraiser"
$ws.Range("D23").Value = 0.4669564962387085
$ws.Range("C24").Value = "warrant"
$ws.Range("D24").Value = 0.4607550501823425
$ws.Range("C25").Value = "unwilling"
$ws.Range("D25").Value = 0.4281245172023773
$ws.Range("C26").Value = "hospitality"
$ws.Range("D26").Value = 0.4055456817150116
$ws.Range("C27").Value = "advantageous"
$ws.Range("D27").Value = 0.4042564332485199
$ws.Range("C28").Value = "scrutinize"
$ws.Range("D28").Value = 0.4034427106380462
$ws.Range("C29").Value = "heavily"
$ws.Range("D29").Value = 0.3929627537727356
$ws.Range("C30").Value = "accelerate"
$ws.Range("D30").Value = 0.3873977363109588
$ws.Range("C31").Value = "marketing"
$ws.Range("D31").Value = 0.3850983381271362
$ws.Range("C32").Value = "bargain"
$ws.Range("D32").Value = 0.4585401713848114
$ws.Range("C33").Value = "fastfood"
$ws.Range("D33").Value = 0.4450857639312744
$ws.Range("C34").Value = "limited"
$ws.Range("D34").Value = 0.443124383687973
$ws.Range("C35").Value = "restrict"
$ws.Range("D35").Value = 0.4423311650753021
$ws.Range("C36").Value = "qualified"
$ws.Range("D36").Value = 0.4328104555606842
$ws.Range("C37").Value = "biofuel"
$ws.Range("D37").Value = 0.4205729365348816
$ws.Range("C38").Value = "trading"
$ws.Range("D38").Value = 0.4141702353954315
$ws.Range("C39").Value = "brand"
$ws.Range("D39").Value = 0.4115100502967834
$ws.Range("D40").Value = 0.4091064035892486
$ws.Range("C41").Value = "relocation"
$ws.Range("D41").Value = 0.3972481787204742
$ws.Range("D42").Value = 0.5313065648078918
$ws.Range("D43").Value = 0.4494546949863434
$ws.Range("C44").Value = "happen"
$ws.Range("D44").Value = 0.4373063445091247
$ws.Range("C45").Value = "expect"
$ws.Range("D45").Value = 0.4262642562389374
$ws.Range("C46").Value = "initiative"
$ws.Range("D46").Value = 0.4171918630599975
$ws.Range("C47").Value = "nothing"
$ws.Range("D47").Value = 0.4106299579143524
$ws.Range("C48").Value = "improvement"
$ws.Range("D48").Value = 0.4088772535324096
$ws.Range("C49").Value = "indication"
$ws.Range("D49").Value = 0.3911709189414978
$ws.Range("C50").Value = "optimistic"
$ws.Range("D50").Value = 0.3872103989124298
$ws.Range("C51").Value = "floorplan"
$ws.Range("D51").Value = 0.3785135149955749
